$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the new columns --------------------------------------
# Current layout:  F=Valor_producto  G=Part_profesional  H=Revisar
# Target layout:   F=Porc_trans  G=Cost_trans  H=Porc_producto  I=Valor_producto
#                  J=Valor_Neto  K=Part_profesional  L=Revisar
# Inserting 4 blank columns at G:J shifts old G->K and old H->L, while F keeps
# its column letter (just gets renamed/repurposed below).
$ws.Columns("G:J").Insert()

# --- 2. Header row ----------------------------------------------------------
$ws.Range("F1").Value = "Porc_trans"
$ws.Range("G1").Value = "Cost_trans"
$ws.Range("H1").Value = "Porc_producto"
$ws.Range("I1").Value = "Valor_producto"
$ws.Range("J1").Value = "Valor_Neto"

# --- 3. Helper to fill in the new per-row figures ---------------------------
function Set-RowFGHIJK {
    param($ws, $row, $F, $G, $H, $I, $J, $K)
    $ws.Cells.Item($row, 6).Value = $F
    $ws.Cells.Item($row, 7).Value = $G
    if ($H -ne $null) {
        $ws.Cells.Item($row, 8).Value = $H
    }
    $ws.Cells.Item($row, 9).Value = $I
    $ws.Cells.Item($row, 10).Value = $J
    $ws.Cells.Item($row, 11).Value = $K
}

Set-RowFGHIJK $ws 2 0.036 1386 0.106 4081 33033 21175
Set-RowFGHIJK $ws 3 0.036 1260 0.1166 4081 29659 19250
Set-RowFGHIJK $ws 4 0.036 1260 0.1166 4081 29659 19250
Set-RowFGHIJK $ws 5 0.036 3240 0.04534444444444444 4081 82679 40500
Set-RowFGHIJK $ws 6 0.036 180 $null 0 4820 4820
Set-RowFGHIJK $ws 7 0.036 1260 0.1166 4081 29659 19250
Set-RowFGHIJK $ws 8 0.036 3240 0.04534444444444444 4081 82679 40500
Set-RowFGHIJK $ws 9 0.036 1260 0.1166 4081 29659 19250
Set-RowFGHIJK $ws 10 0.036 1386 0.106 4081 33033 21175
Set-RowFGHIJK $ws 11 0.036 1008 0.14575 4081 22911 15400
Set-RowFGHIJK $ws 12 0.036 1386 0.106 4081 33033 21175
Set-RowFGHIJK $ws 13 0.036 15120 0.1066452380952381 44791 360089 186209
Set-RowFGHIJK $ws 14 0.036 1620 0.09068888888888889 4081 39299 24750
Set-RowFGHIJK $ws 15 0 0 0.14575 4081 23919 15400
Set-RowFGHIJK $ws 16 0 0 0.0742 4081 50919 30250
Set-RowFGHIJK $ws 17 0 0 0.106 4081 34419 21175
Set-RowFGHIJK $ws 18 0 0 0.14575 4081 23919 15400
Set-RowFGHIJK $ws 19 0.036 1620 0.09068888888888889 4081 39299 24750
Set-RowFGHIJK $ws 20 0.036 1008 0.14575 4081 22911 15400

# --- 4. New discount rows at the bottom -------------------------------------
# The date-like labels in column A ("2024-11-09", "2024-11-15", as used
# already by rows 21-25) must stay plain text, exactly like the existing
# rows above them, rather than being auto-converted to date serials by a
# direct .Value assignment. Route them through a scratch formula cell +
# copy/paste-values so the destination keeps its default (unstyled) text
# cell, just like the rest of the sheet.
function Set-TextValue {
    param($ws, $a1, $text)
    $scratch = $ws.Range("Z1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($a1).PasteSpecial(-4163)
    $scratch.Clear()
}

Set-TextValue $ws "A26" "2024-11-09"
$ws.Range("C26").Value = "Descuento - Producto - Productos de Color"
$ws.Range("D26").Value = "Elvis Molina"
$ws.Range("K26").Value = -50361

Set-TextValue $ws "A27" "2024-11-15"
$ws.Range("C27").Value = "Descuento - Producto - Productos de Color"
$ws.Range("D27").Value = "Elvis Molina"
$ws.Range("K27").Value = -35895
